$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 20000
$ws.Range("I20").Value = 20000
$ws.Range("K20").Value = 20000
$ws.Range("M20").Value = -19770
$ws.Range("H21").Value = 47624.75
$ws.Range("I21").Value = 39249.5
$ws.Range("K21").Value = 39249.5
$ws.Range("M21").Value = -38781.5
$ws.Range("H23").Value = 47624.75
$ws.Range("I23").Value = 39249.5
$ws.Range("K23").Value = 39249.5
$ws.Range("M23").Value = -39015.5
$ws.Range("H34").Value = 18250
$ws.Range("I34").Value = 18250
$ws.Range("K34").Value = 18250
$ws.Range("M34").Value = -18047
$ws.Range("H35").Value = 20000
$ws.Range("I35").Value = 20000
$ws.Range("K35").Value = 20000
$ws.Range("M35").Value = -19621
$ws.Range("H36").Value = 18250
$ws.Range("I36").Value = 18250
$ws.Range("K36").Value = 18250
$ws.Range("M36").Value = -17535

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1668.1177
$ws.Range("J2").Value = 4078.25
$ws.Range("L2").Value = 4078.25
$ws.Range("N2").Value = -4304.25
$ws.Range("H45").Value = 3007.7273
$ws.Range("I45").Value = 1938.3334
$ws.Range("J45").Value = 5299.2856
$ws.Range("K45").Value = 1938.3334
$ws.Range("L45").Value = 5299.2856
$ws.Range("M45").Value = -1561.3334
$ws.Range("N45").Value = -6053.2856
$ws.Range("H61").Value = 9138294
$ws.Range("I61").Value = 10534255
$ws.Range("J61").Value = 2507478.5
$ws.Range("K61").Value = 10534255
$ws.Range("L61").Value = 2507478.5
$ws.Range("M61").Value = -10534043
$ws.Range("N61").Value = -2507902.5
$ws.Range("H74").Value = 2397.7083
$ws.Range("I74").Value = 2352.1162
$ws.Range("K74").Value = 2352.1162
$ws.Range("M74").Value = -1478.1162
$ws.Range("H77").Value = 2397.7083
$ws.Range("I77").Value = 2352.1162
$ws.Range("K77").Value = 11760.581
$ws.Range("M77").Value = -7392.581
$ws.Range("H110").Value = 5758.846
$ws.Range("I110").Value = 5897.15
$ws.Range("K110").Value = 5897.15
$ws.Range("M110").Value = -3852.15
$ws.Range("H116").Value = 1668.1177
$ws.Range("J116").Value = 4078.25
$ws.Range("L116").Value = 4078.25
$ws.Range("N116").Value = -8666.25
$ws.Range("H122").Value = 3004.5715
$ws.Range("I122").Value = 2787.1538
$ws.Range("K122").Value = 8361.4614
$ws.Range("M122").Value = -5911.4614
$ws.Range("H136").Value = 9138294
$ws.Range("I136").Value = 10534255
$ws.Range("J136").Value = 2507478.5
$ws.Range("K136").Value = 31602765
$ws.Range("L136").Value = 7522435.5
$ws.Range("M136").Value = -31600215
$ws.Range("N136").Value = -7527535.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1668.1177
$ws.Range("J3").Value = 4078.25
$ws.Range("L3").Value = 4078.25
$ws.Range("N3").Value = -4306.25
$ws.Range("H95").Value = 41445.2
$ws.Range("J95").Value = 41445.2
$ws.Range("L95").Value = 41445.2
$ws.Range("N95").Value = -46937.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27030700
$ws.Range("J31").Value = 4125
$ws.Range("L31").Value = 4125
$ws.Range("N31").Value = -4715
$ws.Range("H34").Value = 27030700
$ws.Range("J34").Value = 4125
$ws.Range("L34").Value = 4125
$ws.Range("N34").Value = -4529
$ws.Range("H58").Value = 1966.2927
$ws.Range("I58").Value = 2014.6923
$ws.Range("J58").Value = 1882.4
$ws.Range("K58").Value = 2014.6923
$ws.Range("L58").Value = 1882.4
$ws.Range("M58").Value = -1811.6923
$ws.Range("N58").Value = -2288.4
$ws.Range("H94").Value = 1767
$ws.Range("I94").Value = 1662.3334
$ws.Range("K94").Value = 1662.3334
$ws.Range("M94").Value = -1211.3334
$ws.Range("H105").Value = 1221.5
$ws.Range("I105").Value = 1059.8182
$ws.Range("K105").Value = 1059.8182
$ws.Range("M105").Value = 687.1818000000001
$ws.Range("H136").Value = 1966.2927
$ws.Range("I136").Value = 2014.6923
$ws.Range("J136").Value = 1882.4
$ws.Range("K136").Value = 6044.0769
$ws.Range("L136").Value = 5647.200000000001
$ws.Range("M136").Value = -3494.0769
$ws.Range("N136").Value = -10747.2

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 385.63635
$ws.Range("I5").Value = 314.125
$ws.Range("J5").Value = 576.3333
$ws.Range("K5").Value = 942.375
$ws.Range("L5").Value = 1728.9999
$ws.Range("M5").Value = -830.375
$ws.Range("N5").Value = -1952.9999
$ws.Range("H62").Value = 14824.3
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 14824.3
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 44472.89999999999
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -45844.89999999999
$ws.Range("H63").Value = 20614.562
$ws.Range("I63").Value = 13627.25
$ws.Range("K63").Value = 40881.75
$ws.Range("M63").Value = -40132.75
$ws.Range("H65").Value = 14824.3
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 14824.3
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 133418.7
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -140282.7
$ws.Range("H66").Value = 20614.562
$ws.Range("I66").Value = 13627.25
$ws.Range("K66").Value = 122645.25
$ws.Range("M66").Value = -118901.25
$ws.Range("H68").Value = 968.06665
$ws.Range("I68").Value = 759.5
$ws.Range("K68").Value = 2278.5
$ws.Range("M68").Value = -1467.5
$ws.Range("H71").Value = 968.06665
$ws.Range("I71").Value = 759.5
$ws.Range("K71").Value = 6835.5
$ws.Range("M71").Value = -2779.5
$ws.Range("H80").Value = 55559892
$ws.Range("I80").Value = 166666670
$ws.Range("K80").Value = 500000010
$ws.Range("M80").Value = -499999074
$ws.Range("H83").Value = 55559892
$ws.Range("I83").Value = 166666670
$ws.Range("K83").Value = 1500000030
$ws.Range("M83").Value = -1499995350
$ws.Range("J107").Value = 7589472
$ws.Range("L107").Value = 22768416
$ws.Range("N107").Value = -22772256
$ws.Range("H113").Value = 1620.6111
$ws.Range("I113").Value = 1341
$ws.Range("J113").Value = 2060
$ws.Range("K113").Value = 4023
$ws.Range("L113").Value = 6180
$ws.Range("M113").Value = -1853
$ws.Range("N113").Value = -10520
$ws.Range("H135").Value = 385.63635
$ws.Range("I135").Value = 314.125
$ws.Range("J135").Value = 576.3333
$ws.Range("K135").Value = 2827.125
$ws.Range("L135").Value = 5186.9997
$ws.Range("M135").Value = -292.125
$ws.Range("N135").Value = -10256.9997

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 43.692307
$ws.Range("I2").Value = 47
$ws.Range("J2").Value = 36.25
$ws.Range("K2").Value = 47
$ws.Range("L2").Value = 36.25
$ws.Range("M2").Value = 66
$ws.Range("N2").Value = -262.25
$ws.Range("H122").Value = 3816.647
$ws.Range("I122").Value = 3763.3845
$ws.Range("K122").Value = 11290.1535
$ws.Range("M122").Value = -8840.1535
$ws.Range("H131").Value = 87979
$ws.Range("J131").Value = 87979
$ws.Range("L131").Value = 87979
$ws.Range("N131").Value = -98059

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8899.625
$ws.Range("I40").Value = 9439.6
$ws.Range("K40").Value = 9439.6
$ws.Range("M40").Value = -9303.6
$ws.Range("H122").Value = 3797.6875
$ws.Range("I122").Value = 3420.3076
$ws.Range("K122").Value = 10260.9228
$ws.Range("M122").Value = -7810.9228
$ws.Range("H132").Value = 2530.738
$ws.Range("I132").Value = 1345.2903
$ws.Range("K132").Value = 4035.8709
$ws.Range("M132").Value = -1505.8709
$ws.Range("H136").Value = 1958.3055
$ws.Range("I136").Value = 958.62067
$ws.Range("J136").Value = 6099.857
$ws.Range("K136").Value = 2875.86201
$ws.Range("L136").Value = 18299.571
$ws.Range("M136").Value = -325.8620099999998
$ws.Range("N136").Value = -23399.571

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 6973.9614
$ws.Range("I126").Value = 7189.2383
$ws.Range("J126").Value = 6069.8
$ws.Range("K126").Value = 21567.7149
$ws.Range("L126").Value = 18209.4
$ws.Range("M126").Value = -19097.7149
$ws.Range("N126").Value = -23149.4
$ws.Range("H132").Value = 295699.72
$ws.Range("J132").Value = 1668831.5
$ws.Range("L132").Value = 5006494.5
$ws.Range("N132").Value = -5011554.5
$ws.Range("H136").Value = 218127.39
$ws.Range("I136").Value = 712.2683
$ws.Range("K136").Value = 2136.8049
$ws.Range("M136").Value = 413.1950999999999
